# Updated schedule, live coding 8
#
# This script reproduces the OOXML diff: renumbering "Task 17..20" down to
# "Task 15..18" (freeing up slots previously used by a removed Task 19/20
# pair of rows), and adding "live coding" video links to several lessons
# that previously had no video column populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 52: Slide presentations -----------------------------------------
# tasks column (I) changes from the generic "Bonus task" label to "Video",
# and the tasks_link column (M) now points at a new live-coding recording.
$ws.Range("I52").Value = "Video"
$ws.Range("M52").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=09743375-5dbd-435e-ae66-ace000173c94"

# --- Row 54: Checking your work -------------------------------------------
# Adds a video column + link; the old "Task 15" entry in the tasks column
# moves away (renumbered elsewhere), so this cell is cleared out.
$ws.Range("G54").Value = "Video"
$ws.Range("I54").ClearContents()
$ws.Range("K54").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=d0565af8-95d9-495f-a782-ace0001aa7d3"

# --- Row 57: Dynamic graphics ----------------------------------------------
# This lesson previously only had E/D/I/M populated; it now gets full
# notes/video/slides columns + links, while its old "Task 16" entry in the
# tasks column is removed (renumbered elsewhere).
$ws.Range("F57").Value = "Notes"
$ws.Range("G57").Value = "Video"
$ws.Range("H57").Value = "Slides"
$ws.Range("I57").ClearContents()
$ws.Range("J57").Value = "#dynamic"
$ws.Range("K57").Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=8f0f14d9-ddfc-47b5-87fb-ace0002155dd"
$ws.Range("L57").Value = "slides/26-dynamic-graphics.html"

# --- Row 59: Making maps ----------------------------------------------------
# tasks_link bonus-task anchor renumbered (#task-b10 -> #task-b9)
$ws.Range("M59").Value = "#task-b9"

# --- Row 60: "Tasks 15 and 16 due" divider ----------------------------------
# Shifts from week 11 / 2021-05-17 to week 12 / 2021-05-24 (date recomputed
# by the sheet's shared formula once the week number changes).
$ws.Range("A60").Value = 12

# --- Row 62: More about maps ------------------------------------------------
# Task renumbered 17 -> 15
$ws.Range("I62").Value = "Task 15"

# --- Row 64: Alternatives to maps -------------------------------------------
# Task renumbered 18 -> 16
$ws.Range("I64").Value = "Task 16"

# --- Row 65: Factors and Dates ----------------------------------------------
# Task renumbered 19 -> 17
$ws.Range("I65").Value = "Task 17"

# --- Row 66: "Tasks 17 and 18 due" divider ----------------------------------
# Shifts from week 12 / 2021-05-24 to week 13 / 2021-05-31.
$ws.Range("A66").Value = 13

# --- Row 67: Using colour ----------------------------------------------------
# Bonus-task anchor renumbered (#task-b11 -> #task-b10)
$ws.Range("M67").Value = "#task-b10"

# --- Row 70: Themes -----------------------------------------------------------
# Task renumbered 20 -> 18
$ws.Range("I70").Value = "Task 18"

# --- View state: best-effort match of the saved selection ---------------------
$ws.Range("H58").Select()
